$d = $word.ActiveDocument

function Get-Anchor {
    $r = $d.Content
    $r.Find.Execute("Data file (.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Collapse(1)
    return $r
}

# Paragraph 1: leave empty (split point right before "Data file (.")
$r1 = Get-Anchor
$r1.Text = "`r"

# Paragraph 2: new text
$r2 = Get-Anchor
$r2.Text = "Insert some random stuff here and see what happens`r"

# Paragraph 3: empty
$r3 = Get-Anchor
$r3.Text = "`r"

# Paragraph 4: empty
$r4 = Get-Anchor
$r4.Text = "`r"
